# memperbaiki export import manage name
# This script rewrites the "management name list" data rows (rows 2-11,
# columns B:I) so that each row's data is shifted down by one position
# (row 2's original data becomes row 3's data, etc.) and row 2 receives a
# brand new record ("coba"/"cb@gmail.com"/...). Row 11, which previously
# only had partial data, is filled in with what used to be row 10's data,
# and row 10 takes on client_id "080835901" (matching the rest of the
# rows) while keeping its own name/email/gender/department/password data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns B..I across rows 2..11.
# Columns: B=client_id, C=name, D=email, E=gender, F=departement,
#          G=class_id, H=password_prompt, I=password
$rows = @(
    @{ Row = 2;  B = '080835901'; C = 'coba';    D = 'cb@gmail.com';          E = 'Laki-laki'; F = 'pertanian'; G = 5; H = '12345678';                                                           I = '12345678' }
    @{ Row = 3;  B = '080835901'; C = 'test';    D = 'test@gmail.com';        E = 'Perempuan';  F = 'testt';     G = 5; H = 'f5fbc6fe84c365315f491d4275c2f2e5d3c60f25684e1d62e7e9fe63abf8d0d8'; I = 'f5fbc6fe84c365315f491d4275c2f2e5d3c60f25684e1d62e7e9fe63abf8d0d8' }
    @{ Row = 4;  B = '080835901'; C = 'jjl';     D = 'jajal@gmail.com';       E = 'Laki-laki';  F = 'dkv';       G = 5; H = 'jajal1234';                                                          I = 'jajal1234' }
    @{ Row = 5;  B = '080835901'; C = 'jjl';     D = 'jajal@gmail.com';       E = 'Laki-laki';  F = 'dkv';       G = 5; H = 'jajal1234';                                                          I = 'jajal1234' }
    @{ Row = 6;  B = '080835901'; C = 'jjl';     D = 'jajal@gmail.com';       E = 'Laki-laki';  F = 'dkv';       G = 5; H = 'jajal1234';                                                          I = 'jajal1234' }
    @{ Row = 7;  B = '080835901'; C = 'erfer';   D = 'tes@gmail.com';         E = 'laki';       F = 'tbsm';      G = 5; H = 'aveceenaintifadhafirdausming';                                       I = 'aveceenaintifadhafirdausming' }
    @{ Row = 8;  B = '080835901'; C = 'hihi';    D = 'hi@gmail.com';          E = 'laki';       F = 'tkj';       G = 5; H = 'aveceena123';                                                         I = 'aveceena123' }
    @{ Row = 9;  B = '080835901'; C = 'p';       D = 'p@gmail.com';           E = 'Laki';       F = 'tkj';       G = 5; H = '39cc70ddc804005a99b5c9be676f9b550cb3482c57467b05dc88d48857dd7aa3'; I = '39cc70ddc804005a99b5c9be676f9b550cb3482c57467b05dc88d48857dd7aa3' }
    @{ Row = 10; B = '080835901'; C = 'testing'; D = 'test@gmail.com';        E = 'perempuan';  F = 'tkj';       G = 5; H = '8bb0cf6eb9b17d0f7d22b456f121257dc1254e1f01665370476383ea776df414'; I = '8bb0cf6eb9b17d0f7d22b456f121257dc1254e1f01665370476383ea776df414' }
    @{ Row = 11; B = '080835910'; C = 'testing'; D = 'test@gmail.com';        E = 'perempuan';  F = 'tkj';       G = 5; H = '8bb0cf6eb9b17d0f7d22b456f121257dc1254e1f01665370476383ea776df414'; I = '8bb0cf6eb9b17d0f7d22b456f121257dc1254e1f01665370476383ea776df414' }
)

foreach ($r in $rows) {
    $n = $r.Row

    # Columns B, C, D, E, F, H, I hold text data in the source workbook
    # (client_id / name / email / gender / departement / password hashes).
    # Some of those text values look numeric (e.g. "080835901", "12345678")
    # so force the cell format to Text first to avoid Excel auto-converting
    # them to numbers and dropping leading zeros / changing the stored type.
    $textRange = $ws.Range("B$n" + ":F$n")
    $textRange.NumberFormat = "@"
    $ws.Range("H$n" + ":I$n").NumberFormat = "@"

    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
}
